# Apply FRED data refresh update to WTREGEN.xlsx
$wb = $excel.ActiveWorkbook

# --- Update the "Data" sheet: append a new weekly observation row ---
$wsData = $wb.Worksheets.Item("Data")

# Reuse the exact same cell formatting as the preceding observation row
# (A111) so the new row matches the rest of the date column instead of
# creating a new duplicate style entry.
$wsData.Cells.Item(111, 1).Copy() | Out-Null
$wsData.Cells.Item(112, 1).PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# New row 112: date serial 45245 (2023-11-15), value 738.562
$wsData.Cells.Item(112, 1).Value = 45245
$wsData.Cells.Item(112, 2).Value = 738.562

# --- Update the "SeriesInfo" sheet: refresh metadata to match the new pull ---
$wsInfo = $wb.Worksheets.Item("SeriesInfo")

# B3/B4/B7 hold plain "YYYY-MM-DD" text. Assigning such a string straight to
# .Value would be auto-recognized by Excel as a date and stored as a
# formatted serial number instead of text, so instead write it as a text
# formula and immediately collapse it back down to a plain value in place;
# that keeps the cell's type as text (same as the source file) without
# picking up a date number format / style.
$wsInfo.Range("B3").Formula = "=""2023-11-21"""
$wsInfo.Range("B3").Copy() | Out-Null
$wsInfo.Range("B3").PasteSpecial(-4163) | Out-Null  # xlPasteValues

$wsInfo.Range("B4").Formula = "=""2023-11-21"""
$wsInfo.Range("B4").Copy() | Out-Null
$wsInfo.Range("B4").PasteSpecial(-4163) | Out-Null  # xlPasteValues

$wsInfo.Range("B7").Formula = "=""2023-11-15"""
$wsInfo.Range("B7").Copy() | Out-Null
$wsInfo.Range("B7").PasteSpecial(-4163) | Out-Null  # xlPasteValues

# B14 ("last_updated") includes a UTC-offset suffix, so Excel does not treat
# it as a recognizable date/time and a plain .Value assignment keeps it text.
$wsInfo.Range("B14").Value = "2023-11-16 15:34:04-06"

$excel.CutCopyMode = $false
